$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0208 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID cell (B2) from DGS-223 to SCD0011-039
$ws.Range("B2").Value = "SCD0011-039"

# Reset the scrolled view back to the top-left (removes topLeftCell="G1")
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# Update selection to B3
$ws.Range("B3").Select()

# Widen column B so the longer TC_ID value fits (mirrors Excel's best-fit
# recalculation after the cell content changed)
$ws.Columns("B").ColumnWidth = 11.6
